# Updated symbol list on Sat Jan  7 18:40:05 UTC 2023 with GitHub Actions
#
# The crypto price tracker refreshed its scrape: the "Price" (column D) and
# "Volume(1h)" (column E) values for most coin rows moved slightly. Every
# cell in these columns is stored as TEXT (e.g. "261.66", "1.15%"), not as
# a number, so each cell's format is forced to Text ("@") before the new
# value is written - this stops Excel's automatic type inference from
# reinterpreting a numeric-looking string like "261.54" as a number, or a
# string like "1.03%" as a percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (new Price, new Volume(1h)); a $null means "leave this column alone"
$rowUpdates = @(
    @(2,  "261.54",     "1.03%"),
    @(3,  "27.18",      "1.30%"),
    @(4,  "4.700",      "0.64%"),
    @(5,  "0.06189",    "3.23%"),
    @(6,  "6.716",      "0.85%"),
    @(7,  "0.8503",     "-0.74%"),
    @(8,  "0.9164",     "-0.75%"),
    @(9,  "0.1406",     "1.09%"),
    @(10, "0.04649",    "-6.09%"),
    @(11, $null,        "1.04%"),
    @(12, "0.03127",    "2.86%"),
    @(13, "0.09036",    "-1.07%"),
    @(14, "0.001543",   "0.98%"),
    @(15, "0.0006162",  "1.79%"),
    @(16, "0.006070",   "-0.55%"),
    @(17, "3.458",      $null),
    @(18, "3.167",      "0.67%"),
    @(19, "2.196",      "1.28%"),
    @(20, "0.3071",     "-1.22%"),
    @(21, $null,        "0.90%"),
    @(22, "4.095",      "-1.11%"),
    @(23, "0.04239",    "0.11%"),
    @(24, "0.001216",   "0.04%"),
    @(25, $null,        "-5.88%"),
    @(26, $null,        "0.10%"),
    @(27, "0.0001601",  "-6.50%"),
    @(40, "0.03955",    "3.02%"),
    @(41, "0.1112",     "-0.27%"),
    @(42, "0.004123",   "8.37%"),
    @(43, $null,        "-9.70%"),
    @(44, "0.01389",    "-7.79%"),
    @(45, "0.00005136", "0.72%"),
    @(46, $null,        "0.15%"),
    @(48, "0.1677",     "28.72%"),
    @(49, "0.00002101", "0.15%"),
    @(50, "0.0002001",  "0.15%")
)

foreach ($entry in $rowUpdates) {
    $row = $entry[0]
    $price = $entry[1]
    $volume = $entry[2]

    if ($null -ne $price) {
        $cell = $ws.Cells.Item($row, 4)   # column D - Price
        $cell.NumberFormat = "@"
        $cell.Value = $price
    }
    if ($null -ne $volume) {
        $cell = $ws.Cells.Item($row, 5)   # column E - Volume(1h)
        $cell.NumberFormat = "@"
        $cell.Value = $volume
    }
}
